# This script applies the "07-sep" daily update to the EPEX spot price
# workbook:
#   - "Prix Spot" sheet: append a new date column CH (07-sep) with its
#     hourly prices.
#   - "Gaz" sheet: append a new row with the 2025-09-05 gas price.
#   - "CO2" sheet: append a new row with the 2025-09-05 CO2 price.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "Prix Spot" -> new column CH
# ---------------------------------------------------------------------
$wsPrix = $wb.Worksheets.Item("Prix Spot")

# Copy the header formatting (bold, centered, bordered) from the last
# existing date column (CG1) onto the new one (CH1) before writing the
# value, matching how the previous columns are formatted.
$wsPrix.Range("CG1").Copy()
$wsPrix.Range("CH1").PasteSpecial(-4122)
$wsPrix.Range("CH1").Value = "07-sep"

$prixValues = @(
    0,
    -0.01,
    -0.01,
    -0.01,
    -0.01,
    -0.01,
    0,
    5.5,
    0,
    -0.09,
    -4,
    -20.09,
    -40.02,
    -53.4,
    -39.45,
    -19.84,
    -2.4,
    6.08,
    47.5,
    79.02,
    95.09,
    89.06,
    87.68000000000001,
    87.08
)

for ($i = 0; $i -lt $prixValues.Length; $i++) {
    $row = $i + 2
    $wsPrix.Cells.Item($row, 86).Value = $prixValues[$i]
}

# ---------------------------------------------------------------------
# Sheet "Gaz" -> new row 83
# ---------------------------------------------------------------------
$wsGaz = $wb.Worksheets.Item("Gaz")

# Force the date to be stored as plain text (like the rest of column A)
# instead of letting Excel auto-convert the "2025-09-05" string into a
# date serial number.
$wsGaz.Range("A83").NumberFormat = "@"
$wsGaz.Range("A83").Value = "2025-09-05"
$wsGaz.Range("A83").Style = "Normal"
$wsGaz.Range("B83").Value = 31.1

# ---------------------------------------------------------------------
# Sheet "CO2" -> new row 83
# ---------------------------------------------------------------------
$wsCo2 = $wb.Worksheets.Item("CO2")

$wsCo2.Range("A83").NumberFormat = "@"
$wsCo2.Range("A83").Value = "2025-09-05"
$wsCo2.Range("A83").Style = "Normal"
$wsCo2.Range("B83").Value = 75.59
